$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G4:G7").Value = "2017-01-03 07:08:17"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E4:E7").Value = "ht"
$wsZh.Range("H4:H7").Value = "2017-01-03 07:08:04"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E4:E7").Value = "ht"
$wsDe.Range("H4:H7").Value = "2017-01-03 07:08:17"
